$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knight+Player")

# Add people assignments to the "All Sprites" / "All Sounds" / "Scripts" rows
$ws.Range("B68:E68").Value = "Ayden/Zeke"
$ws.Range("B69:E69").Value = "Nick"
$ws.Range("A70").Value = "Scripts"
$ws.Range("B70:E70").Value = "Noah"

# Clear out the old sandbox-planning notes below (rows 71-81, column A)
$ws.Range("A71").Value = ""
$ws.Range("A73:A78").Value = ""
$ws.Range("A80:A81").Value = ""

# Update the saved cell selection to match the author's final cursor position
$ws.Range("D77").Select()

$wb.Save()
